# fix: createDataSet 경로 수정
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = "일교"
$ws.Range("F2").Value = "EGC7017"
$ws.Range("H2").Value = "태권도"
$ws.Range("I2").Value = "이상현"
$ws.Range("J2").Value = 1
$ws.Range("T2").Value = "Tae Kwon Do"

# Row 3 updates
$ws.Range("F3").Value = "GCR7003"
$ws.Range("H3").Value = "일본,문학,그리고불교"
$ws.Range("I3").Value = "김호성(法雨)"
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = "B+"
$ws.Range("T3").Value = "Japan : Literature, and Buddhism"
